# Update GST column (D) handling to flexible search in generate_master_payout.py
# This updates the header label in D1 and the computed GST values in D2:D24
# on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text for column D
$ws.Range("D1").Value = "(19) GST paid by Zomato on behalf of restaurant - under section 9(5)"

# Update the GST values for rows 2-24
$dValues = @{
    2  = 12784.6
    3  = 17158.95
    4  = 3877.05
    5  = 7233.4
    6  = 15055.63
    7  = 19120.08
    8  = 14759.22
    9  = 7412.77
    10 = 15849.58
    11 = 12087.24
    12 = 15803.09
    13 = 12946.84
    14 = 2781.76
    15 = 12998.34
    16 = 20984.57
    17 = 12879.2
    18 = 6852.039999999999
    19 = 23108.32
    20 = 22750.27
    21 = 20797.95
    22 = 12529.67
    23 = 9898.040000000001
    24 = 29601
}

foreach ($row in $dValues.Keys) {
    $ws.Cells.Item($row, 4).Value = $dValues[$row]
}
